$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Add a new Job Posting row with Job_Id = JD_006
$ws.Range("A7").Value = "JD_006"
$ws.Range("B7").Value = "Senior Engineer"
$ws.Range("C7").Value = "dsfdsf"
$ws.Range("D7").Value = 1
$ws.Range("E7").Value = 3
